$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 451, shifting existing rows 451-480 down to 453-482
$ws.Rows.Item(451).Resize(2).Insert()

# Populate new row 451
$ws.Cells.Item(451, 1).Value = 10
$ws.Cells.Item(451, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(451, 3).Value = "La Araucanía"
$ws.Cells.Item(451, 4).Value = 44746
$ws.Cells.Item(451, 5).Value = 9
$ws.Cells.Item(451, 6).Value = 100112032
$ws.Cells.Item(451, 7).Value = "Zapallo italiano"
$ws.Cells.Item(451, 8).Value = "Bola 8"
$ws.Cells.Item(451, 9).Value = "Primera"
$ws.Cells.Item(451, 10).Value = 110
$ws.Cells.Item(451, 11).Value = 10000
$ws.Cells.Item(451, 12).Value = 10000
$ws.Cells.Item(451, 13).Value = 10000
$ws.Cells.Item(451, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(451, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(451, 16).Value = 167
$ws.Cells.Item(451, 17).Value = 60
$ws.Cells.Item(451, 18).Value = "Hortaliza"

# Populate new row 452
$ws.Cells.Item(452, 1).Value = 10
$ws.Cells.Item(452, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(452, 3).Value = "La Araucanía"
$ws.Cells.Item(452, 4).Value = 44746
$ws.Cells.Item(452, 5).Value = 9
$ws.Cells.Item(452, 6).Value = 100112032
$ws.Cells.Item(452, 7).Value = "Zapallo italiano"
$ws.Cells.Item(452, 8).Value = "Sin especificar"
$ws.Cells.Item(452, 9).Value = "Primera"
$ws.Cells.Item(452, 10).Value = 325
$ws.Cells.Item(452, 11).Value = 13000
$ws.Cells.Item(452, 12).Value = 15000
$ws.Cells.Item(452, 13).Value = 13769
$ws.Cells.Item(452, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(452, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(452, 16).Value = 229
$ws.Cells.Item(452, 17).Value = 60
$ws.Cells.Item(452, 18).Value = "Hortaliza"
